$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B3").Value = 8605
$ws.Range("C3").Value = 17200
$ws.Range("D3").Value = 30100
$ws.Range("E3").Value = 36100
$ws.Range("F3").Value = 33000
$ws.Range("G3").Value = 27400
$ws.Range("B4").Value = 1128.267776
$ws.Range("C4").Value = 2255.486976
$ws.Range("D4").Value = 3947.88864
$ws.Range("E4").Value = 4730.126336
$ws.Range("F4").Value = 4455.399424
$ws.Range("G4").Value = 3585.081344
$ws.Range("B5").Value = 115.11
$ws.Range("C5").Value = 115.02
$ws.Range("D5").Value = 131.27
$ws.Range("E5").Value = 206.65
$ws.Range("F5").Value = 455.6
$ws.Range("G5").Value = 1146.99
$ws.Range("B6").Value = 188
$ws.Range("C6").Value = 194
$ws.Range("D6").Value = 241
$ws.Range("E6").Value = 453
$ws.Range("F6").Value = 1037
$ws.Range("G6").Value = 2442
$ws.Range("B7").Value = 190
$ws.Range("C7").Value = 227
$ws.Range("D7").Value = 265
$ws.Range("E7").Value = 635
$ws.Range("F7").Value = 1205
$ws.Range("G7").Value = 2769
$ws.Range("B12").Value = 5038
$ws.Range("C12").Value = 10400
$ws.Range("D12").Value = 17300
$ws.Range("E12").Value = 29900
$ws.Range("F12").Value = 63000
$ws.Range("G12").Value = 389000
$ws.Range("B13").Value = 20.6569472
$ws.Range("C13").Value = 42.6770432
$ws.Range("D13").Value = 70.67402240000001
$ws.Range("E13").Value = 122.683392
$ws.Range("F13").Value = 262.144
$ws.Range("G13").Value = 1592.786944
$ws.Range("B14").Value = 198.16194
$ws.Range("C14").Value = 189.51458
$ws.Range("D14").Value = 223.29757
$ws.Range("E14").Value = 260.96951
$ws.Range("F14").Value = 246.39025
$ws.Range("G14").Value = 79.07133999999999
$ws.Range("B15").Value = 1417.216
$ws.Range("C15").Value = 1433.6
$ws.Range("D15").Value = 1449.984
$ws.Range("E15").Value = 1482.752
$ws.Range("G15").Value = 240.64
$ws.Range("B16").Value = 1531.904
$ws.Range("C16").Value = 1531.904
$ws.Range("D16").Value = 1564.672
$ws.Range("E16").Value = 1613.824
$ws.Range("F16").Value = 1712.128
$ws.Range("G16").Value = 354.304
$ws.Range("B21").Value = 11100
$ws.Range("C21").Value = 15100
$ws.Range("D21").Value = 17100
$ws.Range("E21").Value = 20100
$ws.Range("F21").Value = 21000
$ws.Range("G21").Value = 21500
$ws.Range("B22").Value = 1458.569216
$ws.Range("C22").Value = 1973.420032
$ws.Range("D22").Value = 2236.612608
$ws.Range("E22").Value = 2638.217216
$ws.Range("F22").Value = 2753.560576
$ws.Range("G22").Value = 2812.280832
$ws.Range("B23").Value = 24.87
$ws.Range("C23").Value = 22.51
$ws.Range("D23").Value = 27.75
$ws.Range("E23").Value = 39.14
$ws.Range("F23").Value = 77.91
$ws.Range("G23").Value = 145.16
$ws.Range("B24").Value = 25
$ws.Range("C24").Value = 30
$ws.Range("D24").Value = 36
$ws.Range("E24").Value = 57
$ws.Range("F24").Value = 124
$ws.Range("G24").Value = 180
$ws.Range("B25").Value = 109
$ws.Range("C25").Value = 38
$ws.Range("D25").Value = 73
$ws.Range("E25").Value = 81
$ws.Range("F25").Value = 149
$ws.Range("G25").Value = 392
$ws.Range("B30").Value = 395000
$ws.Range("C30").Value = 546000
$ws.Range("D30").Value = 690000
$ws.Range("E30").Value = 790000
$ws.Range("F30").Value = 586000
$ws.Range("G30").Value = 683000
$ws.Range("B31").Value = 1616.904192
$ws.Range("C31").Value = 2236.612608
$ws.Range("D31").Value = 2825.91232
$ws.Range("E31").Value = 3233.808384
$ws.Range("F31").Value = 2402.287616
$ws.Range("G31").Value = 2797.600768
$ws.Range("B32").Value = 0.96441
$ws.Range("C32").Value = 1.04086
$ws.Range("D32").Value = 1.14305
$ws.Range("E32").Value = 1.4879
$ws.Range("F32").Value = 2.55293
$ws.Range("G32").Value = 4.3003
$ws.Range("B33").Value = 1.608
$ws.Range("C33").Value = 1.768
$ws.Range("D33").Value = 1.8
$ws.Range("E33").Value = 2.448
$ws.Range("F33").Value = 4.256
$ws.Range("G33").Value = 6.176
$ws.Range("B34").Value = 1.848
$ws.Range("C34").Value = 2.128
$ws.Range("D34").Value = 2.224
$ws.Range("E34").Value = 3.152
$ws.Range("F34").Value = 5.664
$ws.Range("G34").Value = 8.640000000000001
$ws.Range("B39").Value = 3696
$ws.Range("C39").Value = 5885
$ws.Range("D39").Value = 8551
$ws.Range("E39").Value = 10000
$ws.Range("F39").Value = 11100
$ws.Range("G39").Value = 9999
$ws.Range("B40").Value = 484.442112
$ws.Range("C40").Value = 771.751936
$ws.Range("D40").Value = 1120.927744
$ws.Range("E40").Value = 1310.72
$ws.Range("F40").Value = 1459.617792
$ws.Range("G40").Value = 1310.72
$ws.Range("B41").Value = 269.95
$ws.Range("C41").Value = 338.63
$ws.Range("D41").Value = 465.52
$ws.Range("E41").Value = 790.37
$ws.Range("F41").Value = 1433.13
$ws.Range("G41").Value = 3191.52
$ws.Range("B42").Value = 1598
$ws.Range("C42").Value = 1598
$ws.Range("D42").Value = 1876
$ws.Range("E42").Value = 2999
$ws.Range("F42").Value = 6259
$ws.Range("G42").Value = 13173
$ws.Range("B43").Value = 1827
$ws.Range("C43").Value = 1696
$ws.Range("D43").Value = 2180
$ws.Range("E43").Value = 4146
$ws.Range("F43").Value = 7504
$ws.Range("G43").Value = 15008
$ws.Range("B48").Value = 122000
$ws.Range("C48").Value = 172000
$ws.Range("D48").Value = 286000
$ws.Range("E48").Value = 287000
$ws.Range("F48").Value = 290000
$ws.Range("G48").Value = 300000
$ws.Range("B49").Value = 501.219328
$ws.Range("C49").Value = 706.740224
$ws.Range("D49").Value = 1169.16224
$ws.Range("E49").Value = 1174.40512
$ws.Range("F49").Value = 1188.036608
$ws.Range("G49").Value = 1227.882496
$ws.Range("B50").Value = 7.94214
$ws.Range("C50").Value = 11.33679
$ws.Range("D50").Value = 13.81346
$ws.Range("E50").Value = 27.46817
$ws.Range("F50").Value = 54.81034
$ws.Range("G50").Value = 106.33994
$ws.Range("C51").Value = 1.752
$ws.Range("D51").Value = 1.704
$ws.Range("E51").Value = 1.656
$ws.Range("F51").Value = 1.656
$ws.Range("G51").Value = 1.656
$ws.Range("B52").Value = 37.632
$ws.Range("C52").Value = 171.008
$ws.Range("D52").Value = 391.168
$ws.Range("E52").Value = 872.448
$ws.Range("F52").Value = 2179.072
$ws.Range("G52").Value = 4145.152
$ws.Range("B57").Value = 8752
$ws.Range("C57").Value = 10000
$ws.Range("D57").Value = 12900
$ws.Range("E57").Value = 12700
$ws.Range("F57").Value = 14000
$ws.Range("G57").Value = 14100
$ws.Range("B58").Value = 1147.142144
$ws.Range("C58").Value = 1435.500544
$ws.Range("D58").Value = 1693.45024
$ws.Range("E58").Value = 1667.23584
$ws.Range("F58").Value = 1837.105152
$ws.Range("G58").Value = 1849.688064
$ws.Range("B59").Value = 21.65344
$ws.Range("C59").Value = 23.42096
$ws.Range("D59").Value = 27.76456
$ws.Range("E59").Value = 38.03683
$ws.Range("F59").Value = 77.91
$ws.Range("G59").Value = 138.94
$ws.Range("B60").Value = 23.424
$ws.Range("C60").Value = 25.472
$ws.Range("D60").Value = 30.592
$ws.Range("E60").Value = 50.944
$ws.Range("F60").Value = 115
$ws.Range("G60").Value = 165
$ws.Range("B61").Value = 24.96
$ws.Range("C61").Value = 27.008
$ws.Range("D61").Value = 33.536
$ws.Range("E61").Value = 54.528
$ws.Range("F61").Value = 127
$ws.Range("G61").Value = 186
$ws.Range("B66").Value = 262000
$ws.Range("C66").Value = 349000
$ws.Range("D66").Value = 400000
$ws.Range("E66").Value = 419000
$ws.Range("F66").Value = 441000
$ws.Range("G66").Value = 436000
$ws.Range("B67").Value = 1073.741824
$ws.Range("C67").Value = 1428.160512
$ws.Range("D67").Value = 1636.827136
$ws.Range("E67").Value = 1717.567488
$ws.Range("F67").Value = 1807.745024
$ws.Range("G67").Value = 1786.773504
$ws.Range("B68").Value = 0.8619600000000001
$ws.Range("C68").Value = 0.8993200000000001
$ws.Range("D68").Value = 0.97435
$ws.Range("E68").Value = 1.25356
$ws.Range("F68").Value = 2.3606
$ws.Range("G68").Value = 4.0742
$ws.Range("B69").Value = 1.272
$ws.Range("C69").Value = 1.4
$ws.Range("D69").Value = 1.544
$ws.Range("E69").Value = 2.064
$ws.Range("F69").Value = 3.984
$ws.Range("G69").Value = 5.664
$ws.Range("B70").Value = 1.512
$ws.Range("C70").Value = 1.688
$ws.Range("D70").Value = 1.928
$ws.Range("E70").Value = 2.864
$ws.Range("F70").Value = 6.24
$ws.Range("G70").Value = 9.664
